$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4349.273
$ws.Range("J17").Value = 4855.125
$ws.Range("L17").Value = 14565.375
$ws.Range("N17").Value = -14901.375
$ws.Range("H33").Value = 1651.2727
$ws.Range("I33").Value = 538.4
$ws.Range("J33").Value = 2578.6667
$ws.Range("K33").Value = 538.4
$ws.Range("L33").Value = 2578.6667
$ws.Range("M33").Value = -309.4
$ws.Range("N33").Value = -3036.6667
$ws.Range("H53").Value = 514.5
$ws.Range("I53").Value = 715
$ws.Range("J53").Value = 233.8
$ws.Range("K53").Value = 715
$ws.Range("L53").Value = 233.8
$ws.Range("M53").Value = -78
$ws.Range("N53").Value = -1507.8
$ws.Range("H88").Value = 1044.6
$ws.Range("I88").Value = 261
$ws.Range("K88").Value = 261
$ws.Range("M88").Value = 145
$ws.Range("H91").Value = 1044.6
$ws.Range("I91").Value = 261
$ws.Range("K91").Value = 261
$ws.Range("M91").Value = 1143
$ws.Range("H115").Value = 285
$ws.Range("I115").Value = 285
$ws.Range("K115").Value = 855
$ws.Range("M115").Value = 712
$ws.Range("H118").Value = 815.8
$ws.Range("I118").Value = 787.8461
$ws.Range("K118").Value = 2363.5383
$ws.Range("M118").Value = -706.5383000000002
$ws.Range("H125").Value = 5983.1665
$ws.Range("I125").Value = 2975
$ws.Range("K125").Value = 26775
$ws.Range("M125").Value = -24315
$ws.Range("H137").Value = 2051.5386
$ws.Range("I137").Value = 1753.3636
$ws.Range("J137").Value = 3691.5
$ws.Range("K137").Value = 5260.0908
$ws.Range("L137").Value = 11074.5
$ws.Range("M137").Value = -2710.0908
$ws.Range("N137").Value = -16174.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1744.1333
$ws.Range("I2").Value = 1705.1538
$ws.Range("J2").Value = 1997.5
$ws.Range("K2").Value = 1705.1538
$ws.Range("L2").Value = 1997.5
$ws.Range("M2").Value = -1592.1538
$ws.Range("N2").Value = -2223.5
$ws.Range("H5").Value = 24.272728
$ws.Range("I5").Value = 24.166666
$ws.Range("J5").Value = 24.4
$ws.Range("K5").Value = 24.166666
$ws.Range("L5").Value = 24.4
$ws.Range("M5").Value = 87.83333400000001
$ws.Range("N5").Value = -248.4
$ws.Range("H110").Value = 2834.25
$ws.Range("I110").Value = 668.5
$ws.Range("K110").Value = 668.5
$ws.Range("M110").Value = 1376.5
$ws.Range("H116").Value = 1744.1333
$ws.Range("I116").Value = 1705.1538
$ws.Range("J116").Value = 1997.5
$ws.Range("K116").Value = 1705.1538
$ws.Range("L116").Value = 1997.5
$ws.Range("M116").Value = 588.8462
$ws.Range("N116").Value = -6585.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1744.1333
$ws.Range("I3").Value = 1705.1538
$ws.Range("J3").Value = 1997.5
$ws.Range("K3").Value = 1705.1538
$ws.Range("L3").Value = 1997.5
$ws.Range("M3").Value = -1591.1538
$ws.Range("N3").Value = -2225.5
$ws.Range("H4").Value = 24.272728
$ws.Range("I4").Value = 24.166666
$ws.Range("J4").Value = 24.4
$ws.Range("K4").Value = 24.166666
$ws.Range("L4").Value = 24.4
$ws.Range("M4").Value = 90.83333400000001
$ws.Range("N4").Value = -254.4
$ws.Range("H105").Value = 1097
$ws.Range("J105").Value = 1232.6666
$ws.Range("L105").Value = 1232.6666
$ws.Range("N105").Value = -4726.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 849.25
$ws.Range("I33").Value = 849.25
$ws.Range("K33").Value = 849.25
$ws.Range("M33").Value = -470.25
$ws.Range("H41").Value = 20000
$ws.Range("I41").Value = 15000
$ws.Range("J41").Value = 21250
$ws.Range("K41").Value = 15000
$ws.Range("L41").Value = 21250
$ws.Range("M41").Value = -14572
$ws.Range("N41").Value = -22106
$ws.Range("H86").Value = 5395
$ws.Range("I86").Value = 5243.75
$ws.Range("K86").Value = 5243.75
$ws.Range("M86").Value = -4120.75
$ws.Range("H89").Value = 5395
$ws.Range("I89").Value = 5243.75
$ws.Range("K89").Value = 26218.75
$ws.Range("M89").Value = -20602.75
$ws.Range("H105").Value = 841.9286
$ws.Range("I105").Value = 849
$ws.Range("J105").Value = 799.5
$ws.Range("K105").Value = 849
$ws.Range("L105").Value = 799.5
$ws.Range("M105").Value = 898
$ws.Range("N105").Value = -4293.5
$ws.Range("H134").Value = 10752.125
$ws.Range("I134").Value = 10145.286
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 30435.858
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -27900.858
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1837.0646
$ws.Range("J4").Value = 1662.4375
$ws.Range("L4").Value = 4987.3125
$ws.Range("N4").Value = -5211.3125
$ws.Range("H6").Value = 7664.875
$ws.Range("I6").Value = 220
$ws.Range("K6").Value = 660
$ws.Range("M6").Value = -547
$ws.Range("H7").Value = 270.2353
$ws.Range("I7").Value = 157.83333
$ws.Range("J7").Value = 540
$ws.Range("K7").Value = 473.49999
$ws.Range("L7").Value = 1620
$ws.Range("M7").Value = -361.49999
$ws.Range("N7").Value = -1844
$ws.Range("H17").Value = 1330.6
$ws.Range("J17").Value = 1650.75
$ws.Range("L17").Value = 4952.25
$ws.Range("N17").Value = -5290.25
$ws.Range("H39").Value = 4000
$ws.Range("J39").Value = 4000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -12588
$ws.Range("H55").Value = 1916.6666
$ws.Range("J55").Value = 4000
$ws.Range("L55").Value = 12000
$ws.Range("N55").Value = -12354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1835.52
$ws.Range("I80").Value = 1066.4
$ws.Range("K80").Value = 1066.4
$ws.Range("M80").Value = -68.40000000000009
$ws.Range("H83").Value = 1835.52
$ws.Range("I83").Value = 1066.4
$ws.Range("K83").Value = 5332
$ws.Range("M83").Value = -340
$ws.Range("H113").Value = 1724.75
$ws.Range("I113").Value = 1724.75
$ws.Range("K113").Value = 1724.75
$ws.Range("M113").Value = 445.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 30676
$ws.Range("I42").Value = 20000
$ws.Range("J42").Value = 36014
$ws.Range("K42").Value = 20000
$ws.Range("L42").Value = 36014
$ws.Range("M42").Value = -19437
$ws.Range("N42").Value = -37140
$ws.Range("H49").Value = 30676
$ws.Range("I49").Value = 20000
$ws.Range("J49").Value = 36014
$ws.Range("K49").Value = 20000
$ws.Range("L49").Value = 36014
$ws.Range("M49").Value = -19853
$ws.Range("N49").Value = -36308
$ws.Range("H100").Value = 1399
$ws.Range("I100").Value = 998.75
$ws.Range("K100").Value = 998.75
$ws.Range("M100").Value = -457.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1161.8125
$ws.Range("I122").Value = 1161.8125
$ws.Range("K122").Value = 3485.4375
$ws.Range("M122").Value = -1035.4375
$ws.Range("H126").Value = 3323.3333
$ws.Range("I126").Value = 3154.1667
$ws.Range("K126").Value = 9462.500100000001
$ws.Range("M126").Value = -6992.500100000001
